{"js": "// Resume edit: add \"Mongoose, SQL, \" to the TECH SKILLS line that lists\n// database technologies, so it reads \"MongoDB, Mongoose, SQL, PostgreSQL\"\n// instead of \"MongoDB, PostgreSQL\".\nconst body = context.document.body;\n\n// \"PostgreSQL\" only occurs once in the document (in the TECH SKILLS\n// database-technologies line), so it is a safe, unique anchor.\nconst matches = body.search(\"PostgreSQL\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length > 0) {\n  // Insert the new text immediately before \"PostgreSQL\" so the line becomes\n  // \"MongoDB, Mongoose, SQL, PostgreSQL\".\n  matches.items[0].insertText(\"Mongoose, SQL, \", Word.InsertLocation.before);\n  await context.sync();\n}\n", "ps1": "# Resume edit: add \"Mongoose, SQL, \" to the TECH SKILLS line that lists\n# database technologies, so it reads \"MongoDB, Mongoose, SQL, PostgreSQL\"\n# instead of \"MongoDB, PostgreSQL\".\n$d = $word.ActiveDocument\n\n# \"PostgreSQL\" only occurs once in the document (in the TECH SKILLS\n# database-technologies line), so it is a safe, unique anchor.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"PostgreSQL\"\n$rng.Find.MatchCase = $true\n$rng.Find.Forward = $true\n$rng.Find.Wrap = 0\n\n$found = $rng.Find.Execute()\n\nif ($found) {\n    # $rng now spans the matched \"PostgreSQL\" text; insert the new\n    # skills right before it so the line reads\n    # \"MongoDB, Mongoose, SQL, PostgreSQL\".\n    $rng.InsertBefore(\"Mongoose, SQL, \")\n}\n"}
